$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $isNumericLooking = $val -match "^[0-9]+$"
    if ($isNumericLooking) {
        $rng.NumberFormat = "@"
        $rng.Value = $val
        $rng.Style = "Normal"
    } else {
        $rng.Value = $val
    }
}

# Row 13
Set-CellText "A13" '05228742831551'
Set-CellText "K13" '2026-02-25T15:09:14+00:00'
Set-CellText "L13" '05228742831551'

# Row 14
Set-CellText "A14" '05228742830266'
Set-CellText "K14" '2026-02-25T15:09:13+00:00'
Set-CellText "L14" '05228742830266'

# Row 21
Set-CellText "A21" '00KNRKVN'
Set-CellText "K21" '2026-02-25T13:05:22+00:00'
Set-CellText "L21" '00KNRKVN'

# Row 22
Set-CellText "A22" '00KNRKVE'
Set-CellText "K22" '2026-02-25T13:05:20+00:00'
Set-CellText "L22" '00KNRKVE'

# Row 23
Set-CellText "A23" '00KNRKVM'
Set-CellText "L23" '00KNRKVM'

# Row 24
Set-CellText "A24" '00KNRKVK'
Set-CellText "K24" '2026-02-25T13:05:21+00:00'
Set-CellText "L24" '00KNRKVK'

# Row 28
Set-CellText "A28" '00KNRKVH'
Set-CellText "E28" '0031009941675'
Set-CellText "I28" '2026-02-25T10:55:09+01:00'
Set-CellText "L28" '00KNRKVH'
Set-CellText "M28" '{"external_order_id": "0031009941675", "sales_office_id": "0303"}'

# Row 29
Set-CellText "A29" '00KNRL56'
Set-CellText "E29" '0031009946885'
Set-CellText "I29" '2026-02-25T13:05:06+01:00'
Set-CellText "L29" '00KNRL56'
Set-CellText "M29" '{"external_order_id": "0031009946885", "sales_office_id": "0303"}'

# Row 35
Set-CellText "A35" 'UTV528451149'
Set-CellText "E35" '0031009941662'
Set-CellText "I35" '2026-02-25T12:50:50+01:00'
Set-CellText "K35" '2026-02-25T16:05:54+00:00'
Set-CellText "L35" 'UTV528451149'
Set-CellText "M35" '{"external_order_id": "0031009941662", "sales_office_id": "0303"}'

# Row 36
Set-CellText "A36" 'UTV551526725'
Set-CellText "E36" '0031009941757'
Set-CellText "I36" '2026-02-25T12:47:17+01:00'
Set-CellText "K36" '2026-02-25T16:05:56+00:00'
Set-CellText "L36" 'UTV551526725'
Set-CellText "M36" '{"external_order_id": "0031009941757", "sales_office_id": "0303"}'

# Row 37
Set-CellText "A37" 'UTV075508160'
Set-CellText "E37" '0031009943162'
Set-CellText "I37" '2026-02-25T10:29:09+01:00'
Set-CellText "L37" 'UTV075508160'
Set-CellText "M37" '{"external_order_id": "0031009943162", "sales_office_id": "0303"}'

# Row 38
Set-CellText "A38" 'UTV346252856'
Set-CellText "E38" '0031009943162'
Set-CellText "I38" '2026-02-25T10:29:36+01:00'
Set-CellText "K38" '2026-02-25T16:05:55+00:00'
Set-CellText "L38" 'UTV346252856'
Set-CellText "M38" '{"external_order_id": "0031009943162", "sales_office_id": "0303"}'

# Row 43
Set-CellText "A43" 'UTV654468660'
Set-CellText "I43" '2026-02-25T16:27:22+01:00'
Set-CellText "L43" 'UTV654468660'

# Row 44
Set-CellText "A44" 'UTV976189216'
Set-CellText "I44" '2026-02-25T16:27:11+01:00'
Set-CellText "L44" 'UTV976189216'

# Row 57
Set-CellText "A57" '1037212543716U'
Set-CellText "E57" '0048006143986'
Set-CellText "I57" '2026-02-25T14:22:55'
Set-CellText "J57" 'GOW'
Set-CellText "K57" '2026-02-25T13:41:58+00:00'
Set-CellText "L57" '1037212543716U'
Set-CellText "M57" '{"external_order_id": "0048006143986", "sales_office_id": "3023"}'

# Row 58
Set-CellText "A58" '1037212162644U'
Set-CellText "E58" '0048006144271'
Set-CellText "I58" '2026-02-25T10:54:37'
Set-CellText "J58" 'TOR'
Set-CellText "K58" '2026-02-25T10:41:55+00:00'
Set-CellText "L58" '1037212162644U'
Set-CellText "M58" '{"external_order_id": "0048006144271", "sales_office_id": "3023"}'

# Row 71
Set-CellText "A71" '1037213855524U'
Set-CellText "L71" '1037213855524U'

# Row 72
Set-CellText "A72" '1037213862350U'
Set-CellText "L72" '1037213862350U'
